$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# The "Enterprises density (per 1000 people)" row (row 13) holds its
# values as text in the shared strings table ("12.1", "3.9", "16").
# Update them to "12.12", "3.93", "16.04" while keeping the cells as
# plain text (not auto-converted to numbers) and preserving their
# original (default) style.
$targets = @(
    @{ Cell = "B13"; Value = "12.12" },
    @{ Cell = "C13"; Value = "3.93" },
    @{ Cell = "D13"; Value = "16.04" }
)

foreach ($t in $targets) {
    $rng = $ws.Range($t.Cell)
    $origStyle = $rng.Style
    # Leading apostrophe forces Excel to store the value as text rather
    # than coercing the numeric-looking string into a number.
    $rng.Value = "'" + $t.Value
    # Restore the original style so no new number-format/quote-prefix
    # styling is left behind on the cell.
    $rng.Style = $origStyle
}
